$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "59.923.18"
$ws.Range("E2").Value = "  +2.75%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.423.71"
$ws.Range("E3").Value = "  +2.54%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.02%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'552.35"
$ws.Range("E5").Value = "  +0.55%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'138.07"
$ws.Range("E6").Value = "  +3.52%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.03%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  +2.71%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  -0.62%  "

# Row 10 - Toncoin
$ws.Range("E10").Value = "  +0.17%  "

# Row 12 - Cardano
$ws.Range("D12").Value = "'0.356"
$ws.Range("E12").Value = "  +0.10%  "

# Row 13 - Avalanche
$ws.Range("D13").Value = "'25.28"
$ws.Range("E13").Value = "  +4.48%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "2.855.70"
$ws.Range("E14").Value = "  +2.58%  "

# Row 15 - WrappedBTC
$ws.Range("D15").Value = "59.864.87"
$ws.Range("E15").Value = "  +2.90%  "

# Row 16 - ShibaInu
$ws.Range("E16").Value = "  +0.83%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.439.32"
$ws.Range("E17").Value = "  +3.41%  "

# Row 18 - Chainlink
$ws.Range("E18").Value = "  +2.59%  "

# Row 19 - Polkadot
$ws.Range("E19").Value = "  +1.54%  "

# Row 20 - BitcoinCash
$ws.Range("D20").Value = "'331.18"
$ws.Range("E20").Value = "  -0.18%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  -3.75%  "

# Row 22 - Dai
$ws.Range("E22").Value = "  -0.03%  "

# Row 23 - Litecoin
$ws.Range("E23").Value = "  +3.62%  "

# Row 24 - Kaspa
$ws.Range("E24").Value = "  +1.04%  "

# Row 25 - InternetComputer(DFINITY)
$ws.Range("E25").Value = "  +5.35%  "

# Row 26 - Binance-PegBSC-USD
$ws.Range("E26").Value = "  +0.18%  "

# Row 27 - Fetch.AI
$ws.Range("E27").Value = "  +3.93%  "

# Row 28 - PEPE
$ws.Range("D28").Value = "0.0₃0779"
$ws.Range("E28").Value = "  +4.69%  "

# Row 29 - PancakeSwap
$ws.Range("E29").Value = "  +0.35%  "

# Row 30 - Monero
$ws.Range("D30").Value = "'170.02"
$ws.Range("E30").Value = "  -0.53%  "

# Row 32 - EthereumClassic
$ws.Range("D32").Value = "'18.69"
$ws.Range("E32").Value = "  +1.32%  "

# Row 33 - SuiNetwork
$ws.Range("E33").Value = "  +1.14%  "

# Row 34 - USDe
$ws.Range("E34").Value = "  -0.01%  "

# Row 35 - ImmutableX
$ws.Range("E35").Value = "  +4.18%  "

# Row 36 - FirstDigitalUSD
$ws.Range("E36").Value = "  -0.03%  "

# Row 37 - NEARProtocol
$ws.Range("D37").Value = "'4.22"
$ws.Range("E37").Value = "  +1.12%  "

# Row 38 - Stacks
$ws.Range("E38").Value = "  +0.40%  "

# Row 39 - OKB
$ws.Range("D39").Value = "'39.64"
$ws.Range("E39").Value = "  -1.81%  "

# Row 40 - PolygonEcosystemToken
$ws.Range("E40").Value = "  -3.84%  "

# Row 41 - Bittensor
$ws.Range("D41").Value = "'314.01"
$ws.Range("E41").Value = "  +9.09%  "

# Row 42 - Filecoin
$ws.Range("E42").Value = "  +0.18%  "

# Row 43 - Aave
$ws.Range("D43").Value = "'139.18"
$ws.Range("E43").Value = "  -0.79%  "

# Row 44 - Stellar
$ws.Range("D44").Value = "'0.0970"
$ws.Range("E44").Value = "  +1.08%  "

# Row 45 - Hedera
$ws.Range("D45").Value = "'0.0521"
$ws.Range("E45").Value = "  +0.90%  "

# Row 46 - InjectiveProtocol
$ws.Range("D46").Value = "'19.53"
$ws.Range("E46").Value = "  +4.73%  "

# Row 47 - Mantle
$ws.Range("E47").Value = "  +2.28%  "

# Row 48 - VeChain
$ws.Range("E48").Value = "  +1.03%  "

# Row 49 - Polygon
$ws.Range("D49").Value = "'0.391"
$ws.Range("E49").Value = "  -8.50%  "

# Row 50 - EnergySwap
$ws.Range("D50").Value = "'17.64"
$ws.Range("E50").Value = "  +0.72%  "
